$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header column F1 with same style as the other headers (B1:E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("F1").Style = $ws.Range("E1").Style

# Fill in the time_taken values for each data row
$ws.Range("F2").Value = "2021-10-05 13:39:49.494326"
$ws.Range("F3").Value = "2021-10-05 13:39:49.494337"
$ws.Range("F4").Value = "2021-10-05 13:39:49.494341"
$ws.Range("F5").Value = "2021-10-05 13:39:49.494345"
